$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "63.448.47"
$ws.Range("E2").Value2 = "  +0.11%  "
$ws.Range("D3").Value2 = "3.077.46"
$ws.Range("E3").Value2 = "  -0.49%  "
$ws.Range("E4").Value2 = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "545.17"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value2 = "  -0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "139.24"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value2 = "  +1.36%  "
$ws.Range("E7").Value2 = "  -0.01%  "
$ws.Range("D8").Value2 = "3.072.28"
$ws.Range("E8").Value2 = "  -0.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.500"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value2 = "  +0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.157"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value2 = "  +0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "6.42"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value2 = "  +2.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.457"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value2 = "  -2.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.0000225"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value2 = "  +3.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "35.00"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value2 = "  -1.47%  "
$ws.Range("D15").Value2 = "3.576.18"
$ws.Range("E15").Value2 = "  -0.51%  "
$ws.Range("D16").Value2 = "63.415.15"
$ws.Range("E16").Value2 = "  +0.03%  "
$ws.Range("E17").Value2 = "  +0.96%  "
$ws.Range("D18").Value2 = "3.071.98"
$ws.Range("E18").Value2 = "  -0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "6.67"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value2 = "  -1.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "476.08"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value2 = "  -2.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "13.49"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value2 = "  -1.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "0.701"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value2 = "  -2.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "7.10"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value2 = "  -2.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "78.68"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value2 = "  -0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "12.24"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value2 = "  -1.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value2 = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "2.72"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value2 = "  -1.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "7.96"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value2 = "  -6.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value2 = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "26.27"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value2 = "  -1.53%  "
$ws.Range("B31").Value2 = "Mantle"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "1.16"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value2 = "  +3.57%  "
$ws.Range("B32").Value2 = "ImmutableX"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "1.90"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value2 = "  -3.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "59.45"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value2 = "  +1.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "2.31"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value2 = "  -7.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "5.52"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value2 = "  +7.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "6.02"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value2 = "  -0.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "490.06"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value2 = "  -4.28%  "
$ws.Range("D38").Value2 = "3.268.35"
$ws.Range("E38").Value2 = "  +3.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.0404"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value2 = "  +0.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.0797"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value2 = "  -0.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.118"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "8.16"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value2 = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "2.61"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value2 = "  -1.47%  "
$ws.Range("E44").Value2 = "  -2.20%  "
$ws.Range("E45").Value2 = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "25.53"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value2 = "  +0.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "124.13"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value2 = "  +2.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "2.03"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value2 = "  -2.34%  "
$ws.Range("D49").Value2 = "0.0₃0530"
$ws.Range("E49").Value2 = "  +4.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.110"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value2 = "  +0.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "2.03"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value2 = "  -0.36%  "
